# Update the "Package" / footprint labels in column AB of the
# "SOD-323" and "TH" worksheets.
#
#   SOD-323 sheet: rows 2-26 "SOD-323"  -> "SMD-SOD323"
#                  row  27   "SOD-123"  -> "SMD-SOD123"
#   TH sheet:      rows 2-3  "DO-41"    -> "TH-DO-41-AND-DO-204AL"
#                  rows 4-7  "DO-201AD" -> "TH-DO-201AD"

$wb = $excel.ActiveWorkbook

$wsSod = $wb.Worksheets.Item("SOD-323")
$wsTh = $wb.Worksheets.Item("TH")

# Order chosen to reproduce the shared-string table ordering of the
# original edit (new strings are appended to the shared string table
# in first-use order):
#   243 TH-DO-41-AND-DO-204AL
#   244 SMD-SOD123
#   245 TH-DO-201AD
#   246 SMD-SOD323
for ($r = 2; $r -le 3; $r++) {
    $wsTh.Range("AB$r").Value = "TH-DO-41-AND-DO-204AL"
}

$wsSod.Range("AB27").Value = "SMD-SOD123"

for ($r = 4; $r -le 7; $r++) {
    $wsTh.Range("AB$r").Value = "TH-DO-201AD"
}

for ($r = 2; $r -le 26; $r++) {
    $wsSod.Range("AB$r").Value = "SMD-SOD323"
}
